$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Swap rows 140 and 141 (and rotate 143/144/145) per updated match ordering ---
# row 140
$ws.Cells.Item(140,1).Value = 138
$ws.Cells.Item(140,2).Value = 7493431
$ws.Cells.Item(140,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(140,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(140,5).Value = 45259.8125
$ws.Cells.Item(140,6).Value = "Sportivo Trinidense"
$ws.Cells.Item(140,7).Value = "Guairena FC"
$ws.Cells.Item(140,8).Value = 7
$ws.Cells.Item(140,9).Value = 2
$ws.Cells.Item(140,10).Value = "H"
$ws.Cells.Item(140,11).Value = 2.05
$ws.Cells.Item(140,12).Value = 3.3
$ws.Cells.Item(140,13).Value = 3.3
$ws.Cells.Item(140,14).Value = 2.6
$ws.Cells.Item(140,15).Value = 3.1
$ws.Cells.Item(140,16).Value = 2.6
$ws.Cells.Item(140,17).Value = 0
$ws.Cells.Item(140,18).Value = 1.925
$ws.Cells.Item(140,19).Value = 1.875
$ws.Cells.Item(140,20).Value = 2.5
$ws.Cells.Item(140,21).Value = 2
$ws.Cells.Item(140,22).Value = 1.8
$ws.Cells.Item(140,23).Value = 1.6
$ws.Cells.Item(140,24).Value = -1
$ws.Cells.Item(140,25).Value = -1
$ws.Cells.Item(140,26).Value = 0.925
$ws.Cells.Item(140,27).Value = -1
$ws.Cells.Item(140,28).Value = 1
$ws.Cells.Item(140,29).Value = -1

# row 141
$ws.Cells.Item(141,1).Value = 139
$ws.Cells.Item(141,2).Value = 7493310
$ws.Cells.Item(141,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(141,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(141,5).Value = 45259.8125
$ws.Cells.Item(141,6).Value = "Libertad Asuncion"
$ws.Cells.Item(141,7).Value = "Tacuary"
$ws.Cells.Item(141,8).Value = 1
$ws.Cells.Item(141,9).Value = 2
$ws.Cells.Item(141,10).Value = "A"
$ws.Cells.Item(141,11).Value = 1.363
$ws.Cells.Item(141,12).Value = 5
$ws.Cells.Item(141,13).Value = 7
$ws.Cells.Item(141,14).Value = 1.571
$ws.Cells.Item(141,15).Value = 4.2
$ws.Cells.Item(141,16).Value = 4.75
$ws.Cells.Item(141,17).Value = -0.75
$ws.Cells.Item(141,18).Value = 1.8
$ws.Cells.Item(141,19).Value = 2
$ws.Cells.Item(141,20).Value = 2.75
$ws.Cells.Item(141,21).Value = 1.8
$ws.Cells.Item(141,22).Value = 2
$ws.Cells.Item(141,23).Value = -1
$ws.Cells.Item(141,24).Value = -1
$ws.Cells.Item(141,25).Value = 3.75
$ws.Cells.Item(141,26).Value = -1
$ws.Cells.Item(141,27).Value = 1
$ws.Cells.Item(141,28).Value = 0.4
$ws.Cells.Item(141,29).Value = -0.5

# row 143
$ws.Cells.Item(143,1).Value = 141
$ws.Cells.Item(143,2).Value = 7493312
$ws.Cells.Item(143,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(143,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(143,5).Value = 45261.8125
$ws.Cells.Item(143,6).Value = "Cerro Porteno"
$ws.Cells.Item(143,7).Value = "Guarani Asuncion"
$ws.Cells.Item(143,8).Value = 4
$ws.Cells.Item(143,9).Value = 0
$ws.Cells.Item(143,10).Value = "H"
$ws.Cells.Item(143,11).Value = 1.7
$ws.Cells.Item(143,12).Value = 3.6
$ws.Cells.Item(143,13).Value = 4.333
$ws.Cells.Item(143,14).Value = 1.727
$ws.Cells.Item(143,15).Value = 3.75
$ws.Cells.Item(143,16).Value = 4.2
$ws.Cells.Item(143,17).Value = -0.5
$ws.Cells.Item(143,18).Value = 1.8
$ws.Cells.Item(143,19).Value = 2
$ws.Cells.Item(143,20).Value = 2.75
$ws.Cells.Item(143,21).Value = 1.875
$ws.Cells.Item(143,22).Value = 1.925
$ws.Cells.Item(143,23).Value = 0.7270000000000001
$ws.Cells.Item(143,24).Value = -1
$ws.Cells.Item(143,25).Value = -1
$ws.Cells.Item(143,26).Value = 0.8
$ws.Cells.Item(143,27).Value = -1
$ws.Cells.Item(143,28).Value = 0.875
$ws.Cells.Item(143,29).Value = -1

# row 144
$ws.Cells.Item(144,1).Value = 142
$ws.Cells.Item(144,2).Value = 7493311
$ws.Cells.Item(144,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(144,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(144,5).Value = 45261.8125
$ws.Cells.Item(144,6).Value = "General Caballero JLM"
$ws.Cells.Item(144,7).Value = "Olimpia Asuncion"
$ws.Cells.Item(144,8).Value = 0
$ws.Cells.Item(144,9).Value = 1
$ws.Cells.Item(144,10).Value = "A"
$ws.Cells.Item(144,11).Value = 3.4
$ws.Cells.Item(144,12).Value = 3.3
$ws.Cells.Item(144,13).Value = 2
$ws.Cells.Item(144,14).Value = 3.2
$ws.Cells.Item(144,15).Value = 3.25
$ws.Cells.Item(144,16).Value = 2.1
$ws.Cells.Item(144,17).Value = 0.25
$ws.Cells.Item(144,18).Value = 1.95
$ws.Cells.Item(144,19).Value = 1.85
$ws.Cells.Item(144,20).Value = 2.25
$ws.Cells.Item(144,21).Value = 1.775
$ws.Cells.Item(144,22).Value = 2.025
$ws.Cells.Item(144,23).Value = -1
$ws.Cells.Item(144,24).Value = -1
$ws.Cells.Item(144,25).Value = 1.1
$ws.Cells.Item(144,26).Value = -1
$ws.Cells.Item(144,27).Value = 0.8500000000000001
$ws.Cells.Item(144,28).Value = -1
$ws.Cells.Item(144,29).Value = 1.025

# row 145
$ws.Cells.Item(145,1).Value = 143
$ws.Cells.Item(145,2).Value = 7493433
$ws.Cells.Item(145,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(145,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(145,5).Value = 45261.8125
$ws.Cells.Item(145,6).Value = "Sportivo Luqueno"
$ws.Cells.Item(145,7).Value = "Nacional Asuncion"
$ws.Cells.Item(145,8).Value = 1
$ws.Cells.Item(145,9).Value = 1
$ws.Cells.Item(145,10).Value = "D"
$ws.Cells.Item(145,11).Value = 2.75
$ws.Cells.Item(145,12).Value = 3.2
$ws.Cells.Item(145,13).Value = 2.4
$ws.Cells.Item(145,14).Value = 2.75
$ws.Cells.Item(145,15).Value = 3.1
$ws.Cells.Item(145,16).Value = 2.45
$ws.Cells.Item(145,17).Value = 0.25
$ws.Cells.Item(145,18).Value = 1.75
$ws.Cells.Item(145,19).Value = 2.05
$ws.Cells.Item(145,20).Value = 2.25
$ws.Cells.Item(145,21).Value = 2
$ws.Cells.Item(145,22).Value = 1.8
$ws.Cells.Item(145,23).Value = -1
$ws.Cells.Item(145,24).Value = 2.1
$ws.Cells.Item(145,25).Value = -1
$ws.Cells.Item(145,26).Value = 0.375
$ws.Cells.Item(145,27).Value = -0.5
$ws.Cells.Item(145,28).Value = -0.5
$ws.Cells.Item(145,29).Value = 0.4

# --- Insert a new match row at 204, shifting subsequent rows down ---
$ws.Rows.Item(204).Insert()
$ws.Range("A203:AC203").Copy()
$ws.Range("A204:AC204").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# row 204
$ws.Cells.Item(204,1).Value = 202
$ws.Cells.Item(204,2).Value = 7609202
$ws.Cells.Item(204,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(204,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(204,5).Value = 45378.79166666666
$ws.Cells.Item(204,6).Value = "Tacuary"
$ws.Cells.Item(204,7).Value = "Sportivo Ameliano"
$ws.Cells.Item(204,8).Value = 1
$ws.Cells.Item(204,9).Value = 1
$ws.Cells.Item(204,10).Value = "D"
$ws.Cells.Item(204,11).Value = 2.75
$ws.Cells.Item(204,12).Value = 3.2
$ws.Cells.Item(204,13).Value = 2.375
$ws.Cells.Item(204,14).Value = 2.625
$ws.Cells.Item(204,15).Value = 3.2
$ws.Cells.Item(204,16).Value = 2.45
$ws.Cells.Item(204,17).Value = 0
$ws.Cells.Item(204,18).Value = 2
$ws.Cells.Item(204,19).Value = 1.8
$ws.Cells.Item(204,20).Value = 2.25
$ws.Cells.Item(204,21).Value = 1.85
$ws.Cells.Item(204,22).Value = 1.95
$ws.Cells.Item(204,23).Value = -1
$ws.Cells.Item(204,24).Value = 2.2
$ws.Cells.Item(204,25).Value = -1
$ws.Cells.Item(204,26).Value = 0
$ws.Cells.Item(204,27).Value = -0
$ws.Cells.Item(204,28).Value = -0.5
$ws.Cells.Item(204,29).Value = 0.475

# row 205
$ws.Cells.Item(205,1).Value = 203
$ws.Cells.Item(205,2).Value = 7609147
$ws.Cells.Item(205,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(205,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(205,5).Value = 45378.89583333334
$ws.Cells.Item(205,6).Value = "Sportivo Trinidense"
$ws.Cells.Item(205,7).Value = "Sportivo Luqueno"
$ws.Cells.Item(205,8).Value = 0
$ws.Cells.Item(205,9).Value = 1
$ws.Cells.Item(205,10).Value = "A"
$ws.Cells.Item(205,11).Value = 3.6
$ws.Cells.Item(205,12).Value = 3.3
$ws.Cells.Item(205,13).Value = 1.909
$ws.Cells.Item(205,14).Value = 3
$ws.Cells.Item(205,15).Value = 3.25
$ws.Cells.Item(205,16).Value = 2.2
$ws.Cells.Item(205,17).Value = 0.25
$ws.Cells.Item(205,18).Value = 1.85
$ws.Cells.Item(205,19).Value = 1.95
$ws.Cells.Item(205,20).Value = 2.5
$ws.Cells.Item(205,21).Value = 1.85
$ws.Cells.Item(205,22).Value = 1.85
$ws.Cells.Item(205,23).Value = -1
$ws.Cells.Item(205,24).Value = -1
$ws.Cells.Item(205,25).Value = 1.2
$ws.Cells.Item(205,26).Value = -1
$ws.Cells.Item(205,27).Value = 0.95
$ws.Cells.Item(205,28).Value = -1
$ws.Cells.Item(205,29).Value = 0.8500000000000001

# row 206
$ws.Cells.Item(206,1).Value = 204
$ws.Cells.Item(206,2).Value = 7609201
$ws.Cells.Item(206,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(206,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(206,5).Value = 45381.79166666666
$ws.Cells.Item(206,6).Value = "Cerro Porteno"
$ws.Cells.Item(206,7).Value = "Nacional Asuncion"
$ws.Cells.Item(206,11).Value = 1.5
$ws.Cells.Item(206,12).Value = 4
$ws.Cells.Item(206,13).Value = 6
$ws.Cells.Item(206,14).Value = 1.5
$ws.Cells.Item(206,15).Value = 4
$ws.Cells.Item(206,16).Value = 6
$ws.Cells.Item(206,17).Value = -1
$ws.Cells.Item(206,18).Value = 1.8
$ws.Cells.Item(206,19).Value = 2
$ws.Cells.Item(206,20).Value = 2.5
$ws.Cells.Item(206,21).Value = 1.95
$ws.Cells.Item(206,22).Value = 1.85
$ws.Cells.Item(206,23).Value = 0
$ws.Cells.Item(206,24).Value = 0
$ws.Cells.Item(206,25).Value = 0
$ws.Cells.Item(206,26).Value = 0
$ws.Cells.Item(206,27).Value = 0

# row 207
$ws.Cells.Item(207,1).Value = 205
$ws.Cells.Item(207,2).Value = 7609146
$ws.Cells.Item(207,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(207,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(207,5).Value = 45381.89583333334
$ws.Cells.Item(207,6).Value = "Libertad Asuncion"
$ws.Cells.Item(207,7).Value = "Olimpia Asuncion"
$ws.Cells.Item(207,11).Value = 1.8
$ws.Cells.Item(207,12).Value = 3.3
$ws.Cells.Item(207,13).Value = 4
$ws.Cells.Item(207,14).Value = 1.8
$ws.Cells.Item(207,15).Value = 3.3
$ws.Cells.Item(207,16).Value = 4
$ws.Cells.Item(207,17).Value = -0.5
$ws.Cells.Item(207,18).Value = 1.85
$ws.Cells.Item(207,19).Value = 1.95
$ws.Cells.Item(207,20).Value = 2.25
$ws.Cells.Item(207,21).Value = 1.825
$ws.Cells.Item(207,22).Value = 1.975
$ws.Cells.Item(207,23).Value = 0
$ws.Cells.Item(207,24).Value = 0
$ws.Cells.Item(207,25).Value = 0
$ws.Cells.Item(207,26).Value = 0
$ws.Cells.Item(207,27).Value = 0

# row 208
$ws.Cells.Item(208,1).Value = 206
$ws.Cells.Item(208,2).Value = 7609145
$ws.Cells.Item(208,3).Value = "Paraguay Division Profesional"
$ws.Cells.Item(208,4).Value = "Paraguay Division Profesional"
$ws.Cells.Item(208,5).Value = 45382.77083333334
$ws.Cells.Item(208,6).Value = "Sol de America"
$ws.Cells.Item(208,7).Value = "General Caballero JLM"
$ws.Cells.Item(208,11).Value = 2.1
$ws.Cells.Item(208,12).Value = 3.2
$ws.Cells.Item(208,13).Value = 3.2
$ws.Cells.Item(208,14).Value = 2.6
$ws.Cells.Item(208,15).Value = 3.1
$ws.Cells.Item(208,16).Value = 2.6
$ws.Cells.Item(208,17).Value = 0
$ws.Cells.Item(208,18).Value = 1.9
$ws.Cells.Item(208,19).Value = 1.9
$ws.Cells.Item(208,20).Value = 2.25
$ws.Cells.Item(208,21).Value = 1.9
$ws.Cells.Item(208,22).Value = 1.9
$ws.Cells.Item(208,23).Value = 0
$ws.Cells.Item(208,24).Value = 0
$ws.Cells.Item(208,25).Value = 0
$ws.Cells.Item(208,26).Value = 0
$ws.Cells.Item(208,27).Value = 0

